# "add area to Q files stn3"
#
# Adds cross-sectional Area columns alongside the existing discharge (Q)
# columns on Sheet1:
#   G  -> "Area"   per-segment area            (D.-D.)*B./100
#   H  -> "Atotal" sum of the segment areas     SUM(G2:G11)
#   J  -> "Atotal" copy/echo of the area total  =H2
#   K  -> "Qtotal" copy/echo of the discharge total =F2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# Row 2: first segment area (measured from 0), the two totals, and the
# echoed Qtotal
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# Row 3: next segment area (standalone formula, not part of the fill below)
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# Rows 4-15: remaining segment areas, filled as one block so Excel keeps
# them as a shared formula group (matches the rest of the sheet's style)
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# Leave the active selection on the new header cell, like the source edit
$ws.Range("G1").Select()
